$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(63, 8).Value = 30000
$ws.Cells.Item(63, 10).Value = 30000
$ws.Cells.Item(63, 12).Value = 30000
$ws.Cells.Item(63, 14).Value = -31248
$ws.Cells.Item(66, 8).Value = 30000
$ws.Cells.Item(66, 10).Value = 30000
$ws.Cells.Item(66, 12).Value = 90000
$ws.Cells.Item(66, 14).Value = -96240
$ws.Cells.Item(88, 8).Value = 439486.88
$ws.Cells.Item(88, 9).Value = 973195.25
$ws.Cells.Item(88, 10).Value = 12520.2
$ws.Cells.Item(88, 11).Value = 973195.25
$ws.Cells.Item(88, 12).Value = 12520.2
$ws.Cells.Item(88, 13).Value = -972789.25
$ws.Cells.Item(88, 14).Value = -13332.2
$ws.Cells.Item(91, 8).Value = 439486.88
$ws.Cells.Item(91, 9).Value = 973195.25
$ws.Cells.Item(91, 10).Value = 12520.2
$ws.Cells.Item(91, 11).Value = 973195.25
$ws.Cells.Item(91, 12).Value = 12520.2
$ws.Cells.Item(91, 13).Value = -971791.25
$ws.Cells.Item(91, 14).Value = -15328.2
$ws.Cells.Item(112, 8).Value = 1483.0278
$ws.Cells.Item(112, 10).Value = 1703.6666
$ws.Cells.Item(112, 12).Value = 5110.9998
$ws.Cells.Item(112, 14).Value = -7326.9998
$ws.Cells.Item(132, 8).Value = 23134.842
$ws.Cells.Item(132, 9).Value = 3583.1667
$ws.Cells.Item(132, 10).Value = 65031.285
$ws.Cells.Item(132, 11).Value = 10749.5001
$ws.Cells.Item(132, 12).Value = 195093.855
$ws.Cells.Item(132, 13).Value = -8219.500100000001
$ws.Cells.Item(132, 14).Value = -200153.855
$ws.Cells.Item(137, 8).Value = 1331303.2
$ws.Cells.Item(137, 9).Value = 2408764
$ws.Cells.Item(137, 10).Value = 5197.923
$ws.Cells.Item(137, 11).Value = 7226292
$ws.Cells.Item(137, 12).Value = 15593.769
$ws.Cells.Item(137, 13).Value = -7223742
$ws.Cells.Item(137, 14).Value = -20693.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10601.528
$ws.Cells.Item(32, 9).Value = 9891
$ws.Cells.Item(32, 10).Value = 15417.333
$ws.Cells.Item(32, 11).Value = 9891
$ws.Cells.Item(32, 12).Value = 15417.333
$ws.Cells.Item(32, 13).Value = -9604
$ws.Cells.Item(32, 14).Value = -15991.333
$ws.Cells.Item(45, 8).Value = 1971.7
$ws.Cells.Item(45, 9).Value = 1727.0667
$ws.Cells.Item(45, 10).Value = 2705.6
$ws.Cells.Item(45, 11).Value = 1727.0667
$ws.Cells.Item(45, 12).Value = 2705.6
$ws.Cells.Item(45, 13).Value = -1350.0667
$ws.Cells.Item(45, 14).Value = -3459.6
$ws.Cells.Item(61, 8).Value = 2221.7932
$ws.Cells.Item(61, 9).Value = 1413.5
$ws.Cells.Item(61, 11).Value = 1413.5
$ws.Cells.Item(61, 13).Value = -1201.5
$ws.Cells.Item(63, 8).Value = 3132.25
$ws.Cells.Item(63, 9).Value = 2317
$ws.Cells.Item(63, 10).Value = 4491
$ws.Cells.Item(63, 11).Value = 2317
$ws.Cells.Item(63, 12).Value = 4491
$ws.Cells.Item(63, 13).Value = -1631
$ws.Cells.Item(63, 14).Value = -5863
$ws.Cells.Item(66, 8).Value = 3132.25
$ws.Cells.Item(66, 9).Value = 2317
$ws.Cells.Item(66, 10).Value = 4491
$ws.Cells.Item(66, 11).Value = 11585
$ws.Cells.Item(66, 12).Value = 22455
$ws.Cells.Item(66, 13).Value = -8153
$ws.Cells.Item(66, 14).Value = -29319
$ws.Cells.Item(123, 8).Value = 37714.5
$ws.Cells.Item(123, 10).Value = 37714.5
$ws.Cells.Item(123, 12).Value = 37714.5
$ws.Cells.Item(123, 14).Value = -47514.5
$ws.Cells.Item(136, 8).Value = 2221.7932
$ws.Cells.Item(136, 9).Value = 1413.5
$ws.Cells.Item(136, 11).Value = 4240.5
$ws.Cells.Item(136, 13).Value = -1690.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5808.077
$ws.Cells.Item(20, 9).Value = 1700.1666
$ws.Cells.Item(20, 10).Value = 9329.143
$ws.Cells.Item(20, 11).Value = 1700.1666
$ws.Cells.Item(20, 12).Value = 9329.143
$ws.Cells.Item(20, 13).Value = -1453.1666
$ws.Cells.Item(20, 14).Value = -9823.143
$ws.Cells.Item(86, 8).Value = 2383.125
$ws.Cells.Item(86, 9).Value = 1961
$ws.Cells.Item(86, 11).Value = 1961
$ws.Cells.Item(86, 13).Value = -838
$ws.Cells.Item(89, 8).Value = 2383.125
$ws.Cells.Item(89, 9).Value = 1961
$ws.Cells.Item(89, 11).Value = 9805
$ws.Cells.Item(89, 13).Value = -4189
$ws.Cells.Item(105, 8).Value = 3854.4119
$ws.Cells.Item(105, 9).Value = 3641.8572
$ws.Cells.Item(105, 10).Value = 4003.2
$ws.Cells.Item(105, 11).Value = 3641.8572
$ws.Cells.Item(105, 12).Value = 4003.2
$ws.Cells.Item(105, 13).Value = -1894.8572
$ws.Cells.Item(105, 14).Value = -7497.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3587437.2
$ws.Cells.Item(31, 9).Value = 1644.275
$ws.Cells.Item(31, 10).Value = 6293696
$ws.Cells.Item(31, 11).Value = 1644.275
$ws.Cells.Item(31, 12).Value = 6293696
$ws.Cells.Item(31, 13).Value = -1349.275
$ws.Cells.Item(31, 14).Value = -6294286
$ws.Cells.Item(34, 8).Value = 3587437.2
$ws.Cells.Item(34, 9).Value = 1644.275
$ws.Cells.Item(34, 10).Value = 6293696
$ws.Cells.Item(34, 11).Value = 1644.275
$ws.Cells.Item(34, 12).Value = 6293696
$ws.Cells.Item(34, 13).Value = -1442.275
$ws.Cells.Item(34, 14).Value = -6294100
$ws.Cells.Item(99, 8).Value = 2576.6667
$ws.Cells.Item(99, 9).Value = 2722.8
$ws.Cells.Item(99, 10).Value = 2394
$ws.Cells.Item(99, 11).Value = 2722.8
$ws.Cells.Item(99, 12).Value = 2394
$ws.Cells.Item(99, 13).Value = -1224.8
$ws.Cells.Item(99, 14).Value = -5390
$ws.Cells.Item(122, 8).Value = 76053.44
$ws.Cells.Item(122, 9).Value = 93410.766
$ws.Cells.Item(122, 10).Value = 838.3333
$ws.Cells.Item(122, 11).Value = 280232.298
$ws.Cells.Item(122, 12).Value = 2514.9999
$ws.Cells.Item(122, 13).Value = -277782.298
$ws.Cells.Item(122, 14).Value = -7414.9999
$ws.Cells.Item(126, 8).Value = 2576.6667
$ws.Cells.Item(126, 9).Value = 2722.8
$ws.Cells.Item(126, 10).Value = 2394
$ws.Cells.Item(126, 11).Value = 8168.400000000001
$ws.Cells.Item(126, 12).Value = 7182
$ws.Cells.Item(126, 13).Value = -5698.400000000001
$ws.Cells.Item(126, 14).Value = -12122
$ws.Cells.Item(132, 8).Value = 670315.0600000001
$ws.Cells.Item(132, 9).Value = 2125.5833
$ws.Cells.Item(132, 10).Value = 1561234.4
$ws.Cells.Item(132, 11).Value = 6376.749899999999
$ws.Cells.Item(132, 12).Value = 4683703.199999999
$ws.Cells.Item(132, 13).Value = -3846.749899999999
$ws.Cells.Item(132, 14).Value = -4688763.199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1629.3334
$ws.Cells.Item(34, 9).Value = 525
$ws.Cells.Item(34, 10).Value = 1850.2
$ws.Cells.Item(34, 11).Value = 1575
$ws.Cells.Item(34, 12).Value = 5550.6
$ws.Cells.Item(34, 13).Value = -1491
$ws.Cells.Item(34, 14).Value = -5718.6
$ws.Cells.Item(68, 8).Value = 1377.439
$ws.Cells.Item(68, 9).Value = 1043.579
$ws.Cells.Item(68, 10).Value = 1478.127
$ws.Cells.Item(68, 11).Value = 3130.737
$ws.Cells.Item(68, 12).Value = 4434.380999999999
$ws.Cells.Item(68, 13).Value = -2319.737
$ws.Cells.Item(68, 14).Value = -6056.380999999999
$ws.Cells.Item(71, 8).Value = 1377.439
$ws.Cells.Item(71, 9).Value = 1043.579
$ws.Cells.Item(71, 10).Value = 1478.127
$ws.Cells.Item(71, 11).Value = 9392.210999999999
$ws.Cells.Item(71, 12).Value = 13303.143
$ws.Cells.Item(71, 13).Value = -5336.210999999999
$ws.Cells.Item(71, 14).Value = -21415.143
$ws.Cells.Item(107, 8).Value = 9783.044
$ws.Cells.Item(107, 10).Value = 10027.583
$ws.Cells.Item(107, 12).Value = 30082.749
$ws.Cells.Item(107, 14).Value = -33922.749
$ws.Cells.Item(140, 8).Value = 85354.25
$ws.Cells.Item(140, 9).Value = 126106.375
$ws.Cells.Item(140, 10).Value = 3850
$ws.Cells.Item(140, 11).Value = 378319.125
$ws.Cells.Item(140, 12).Value = 11550
$ws.Cells.Item(140, 13).Value = -373139.125
$ws.Cells.Item(140, 14).Value = -21910

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5241.25
$ws.Cells.Item(70, 9).Value = 5202.9487
$ws.Cells.Item(70, 10).Value = 5540
$ws.Cells.Item(70, 11).Value = 5202.9487
$ws.Cells.Item(70, 12).Value = 5540
$ws.Cells.Item(70, 13).Value = -4932.9487
$ws.Cells.Item(70, 14).Value = -6080
$ws.Cells.Item(73, 8).Value = 5241.25
$ws.Cells.Item(73, 9).Value = 5202.9487
$ws.Cells.Item(73, 10).Value = 5540
$ws.Cells.Item(73, 11).Value = 5202.9487
$ws.Cells.Item(73, 12).Value = 5540
$ws.Cells.Item(73, 13).Value = -4266.9487
$ws.Cells.Item(73, 14).Value = -7412
$ws.Cells.Item(122, 8).Value = 1487.125
$ws.Cells.Item(122, 9).Value = 1316.1666
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 3948.4998
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -1498.4998
$ws.Cells.Item(122, 14).Value = -10900

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 40429
$ws.Cells.Item(123, 10).Value = 40429
$ws.Cells.Item(123, 12).Value = 40429
$ws.Cells.Item(123, 14).Value = -50229
$ws.Cells.Item(132, 8).Value = 2176422.2
$ws.Cells.Item(132, 9).Value = 3346073.2
$ws.Cells.Item(132, 11).Value = 10038219.6
$ws.Cells.Item(132, 13).Value = -10035689.6
$ws.Cells.Item(136, 8).Value = 371202.28
$ws.Cells.Item(136, 9).Value = 519096.8
$ws.Cells.Item(136, 10).Value = 1465.9445
$ws.Cells.Item(136, 11).Value = 1557290.4
$ws.Cells.Item(136, 12).Value = 4397.833500000001
$ws.Cells.Item(136, 13).Value = -1554740.4
$ws.Cells.Item(136, 14).Value = -9497.833500000001
